$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_cal_temp - single calibration temperature value
$ws.Range("H6").Value = 20.14

# CC_wl - wavelength array
$ws.Range("H7").Value = "[188.96, 189.76, 190.56, 191.35, 192.15, 192.95, 193.75, 194.55, 195.35, 196.15, 196.95, 197.76, 198.56, 199.36, 200.16, 200.96, 201.76, 202.56, 203.37, 204.17, 204.97, 205.78, 206.58, 207.38, 208.18, 208.99, 209.79, 210.6, 211.4, 212.2, 213.01, 213.81, 214.62, 215.42, 216.23, 217.03, 217.84, 218.64, 219.45, 220.26, 221.06, 221.87, 222.67, 223.48, 224.29, 225.09, 225.9, 226.71, 227.52, 228.32, 229.13, 229.94, 230.75, 231.55, 232.36, 233.17, 233.98, 234.79, 235.59, 236.4, 237.21, 238.02, 238.83, 239.64, 240.45, 241.26, 242.07, 242.88, 243.69, 244.49, 245.3, 246.11, 246.92, 247.73, 248.54, 249.35, 250.16, 250.98, 251.79, 252.6, 253.41, 254.22, 255.03, 255.84, 256.65, 257.46, 258.27, 259.08, 259.89, 260.71, 261.52, 262.33, 263.14, 263.95, 264.76, 265.57, 266.39, 267.2, 268.01, 268.82, 269.63, 270.44, 271.26, 272.07, 272.88, 273.69, 274.5, 275.32, 276.13, 276.94, 277.75, 278.56, 279.38, 280.19, 281, 281.81, 282.62, 283.44, 284.25, 285.06, 285.87, 286.68, 287.5, 288.31, 289.12, 289.93, 290.75, 291.56, 292.37, 293.18, 293.99, 294.81, 295.62, 296.43, 297.24, 298.05, 298.87, 299.68, 300.49, 301.3, 302.11, 302.92, 303.74, 304.55, 305.36, 306.17, 306.98, 307.79, 308.61, 309.42, 310.23, 311.04, 311.85, 312.66, 313.47, 314.28, 315.1, 315.91, 316.72, 317.53, 318.34, 319.15, 319.96, 320.77, 321.58, 322.39, 323.2, 324.01, 324.82, 325.63, 326.44, 327.25, 328.06, 328.87, 329.68, 330.49, 331.3, 332.11, 332.92, 333.73, 334.54, 335.34, 336.15, 336.96, 337.77, 338.58, 339.39, 340.19, 341, 341.81, 342.62, 343.43, 344.23, 345.04, 345.85, 346.65, 347.46, 348.27, 349.07, 349.88, 350.69, 351.49, 352.3, 353.1, 353.91, 354.72, 355.52, 356.33, 357.13, 357.94, 358.74, 359.55, 360.35, 361.15, 361.96, 362.76, 363.57, 364.37, 365.17, 365.98, 366.78, 367.58, 368.38, 369.19, 369.99, 370.79, 371.59, 372.39, 373.2, 374, 374.8, 375.6, 376.4, 377.2, 378, 378.8, 379.6, 380.4, 381.2, 382, 382.8, 383.6, 384.39, 385.19, 385.99, 386.79, 387.59, 388.38, 389.18, 389.98, 390.77, 391.57, 392.37, 393.16, 393.96, 394.75]"

# CC_eno3 array
$ws.Range("H8").Value = "[0.01192793, -0.00697043, 0.005849, 0.00665215, 0.00751971, -0.00075661, -0.00000902, 0.01191753, 0.0164893, 0.00845698, 0.00751078, -0.00754223, -0.00358466, -0.01460887, -0.01400365, -0.00486794, -0.00529526, 0.009174, 0.01828292, 0.02422251, 0.04105642, 0.00239134, 0.0038067, 0.00560385, 0.00578101, 0.00727442, 0.00733214, 0.00733089, 0.00714702, 0.00685463, 0.0065175, 0.00615574, 0.00579062, 0.00545598, 0.00509853, 0.00474432, 0.00440686, 0.0040569, 0.00372733, 0.00339799, 0.00309102, 0.00278452, 0.00249647, 0.00222409, 0.00196811, 0.00173252, 0.00152083, 0.00131634, 0.00113886, 0.00098552, 0.00084869, 0.00072231, 0.00061484, 0.0005237, 0.00044217, 0.00036545, 0.00029558, 0.00024206, 0.00020124, 0.0001613, 0.00013333, 0.00009658, 0.00007796, 0.00005726, 0.00003693, 0.00002848, 0.00001458, 0.0000094, -0.00000289, -0.00000882, -0.00001219, -0.00001576, -0.00001243, -0.00000779, -0.000017, -0.00001705, -0.00001857, -0.00001012, -0.00001247, -0.00001645, -0.00001141, -0.00000411, 0.00000252, 0.00000146, 0.00001989, 0.00000655, 0.00000867, 0.00001806, 0.00002053, 0.00002164, 0.00001452, 0.0000234, 0.00002651, 0.00003163, 0.00002424, 0.00002865, 0.00003029, 0.00003239, 0.00003755, 0.00003289, 0.00003643, 0.00003723, 0.00004751, 0.0000407, 0.00004822, 0.00003673, 0.00005524, 0.00004783, 0.00005253, 0.00005513, 0.00006347, 0.00005787, 0.00006424, 0.00005386, 0.00007029, 0.00007547, 0.00008689, 0.00008656, 0.00009828, 0.00010867, 0.00010115, 0.00011592, 0.00011065, 0.00011561, 0.00011103, 0.00010706, 0.00012003, 0.00012983, 0.00013147, 0.00014247, 0.00013999, 0.00015111, 0.00014573, 0.00015223, 0.00016555, 0.00016661, 0.00015769, 0.00016258, 0.00016542, 0.0001607, 0.00016966, 0.00017492, 0.00017167, 0.00017214, 0.00017353, 0.0001787, 0.00018075, 0.00018474, 0.0001807, 0.00018247, 0.00018456, 0.00017997, 0.00018131, 0.00018507, 0.00018562, 0.00018731, 0.00018364, 0.00019096, 0.00019405, 0.00020303, 0.00020686, 0.00021817, 0.00020382, 0.00020623, 0.00021363, 0.00021981, 0.00021433, 0.00022896, 0.00022379, 0.00022696, 0.00024402, 0.00022867, 0.00023002, 0.00024863, 0.00023345, 0.00024215, 0.00025062, 0.00025227, 0.00024991, 0.00024461, 0.0002642, 0.00026664, 0.00026795, 0.00027281, 0.00027795, 0.00027396, 0.00028254, 0.00027331, 0.0002854, 0.00028358, 0.00029118, 0.00028896, 0.00029285, 0.00028586, 0.00027106, 0.00030821, 0.0002925, 0.00029893, 0.00029926, 0.00029815, 0.00030348, 0.00028146, 0.00029786, 0.00027699, 0.00029182, 0.00028747, 0.00029992, 0.00029465, 0.00028562, 0.00029514, 0.00029609, 0.00028005, 0.00027774, 0.00026771, 0.00027884, 0.00028807, 0.00028284, 0.00031331, 0.00030071, 0.00030251, 0.00028165, 0.00030224, 0.00029043, 0.00029314, 0.00030383, 0.00031174, 0.00027438, 0.00030111, 0.00026839, 0.00031314, 0.00026808, 0.00026913, 0.00030803, 0.00030828, 0.00029758, 0.00029112, 0.00027748, 0.00028845, 0.00031492, 0.00031715, 0.00030247, 0.00031411, 0.00030916, 0.00035733, 0.00032358, 0.00030727, 0.00027717, 0.00026906, 0.0002818, 0.00025523, 0.00027945, 0.00032143, 0.00028039, 0.00028942, 0.00027121, 0.00029941]"

# CC_eswa array
$ws.Range("H9").Value = "[0.00291964, 0.01317798, 0.00420499, -0.00278269, 0.00033558, 0.00558209, -0.00671788, -0.01152932, -0.00361449, -0.00843627, -0.00339183, 0.01257221, 0.01049906, 0.04424796, 0.05953843, 0.06544391, 0.07611656, 0.07958872, 0.08056865, 0.08262449, 0.07818528, 0.08660406, 0.07907087, 0.07295199, 0.06315188, 0.05310584, 0.04441012, 0.0368419, 0.03034277, 0.02488169, 0.02023406, 0.01639394, 0.01323756, 0.01059064, 0.00845701, 0.00670867, 0.0052945, 0.00415032, 0.00324543, 0.00252205, 0.00195424, 0.00151912, 0.00117859, 0.00091919, 0.00071957, 0.00056862, 0.00045587, 0.00036979, 0.00031087, 0.00024505, 0.00020171, 0.00016479, 0.00014009, 0.00010945, 0.00008908, 0.00007529, 0.00006863, 0.00005796, 0.00005048, 0.00003798, 0.00002848, 0.00002677, 0.00002915, 0.00001397, 0.00001569, 0.00001122, 0.00000688, 0.00000488, -0.00000111, -0.00000234, -0.00000528, -0.0000098, -0.0000052, -0.00000645, -0.00000647, -0.0000124, -0.00000724, -0.00000844, 0.00000072, 0.00000056, 0.00000179, 0.0000057, -0.00000536, 0.00000342, -0.0000073, 0.00000233, 0.00000627, 0.0000104, 0.00000753, 0.00000152, 0.00001854, 0.00001492, 0.00000781, 0.00001729, 0.00002081, 0.0000221, 0.00002596, 0.00002153, 0.00002408, 0.00003171, 0.00002891, 0.00002873, 0.00002092, 0.00002362, 0.00003339, 0.00004807, 0.00004051, 0.00005368, 0.00004767, 0.00005288, 0.00004529, 0.00003816, 0.00004382, 0.00005164, 0.00004163, 0.00005537, 0.00003703, 0.00004621, 0.00003677, 0.00004675, 0.000049, 0.00003758, 0.00006201, 0.00005325, 0.00006239, 0.00007911, 0.00006152, 0.0000605, 0.00006885, 0.00005821, 0.00006785, 0.00006551, 0.00007656, 0.00007561, 0.00007172, 0.00007336, 0.0000903, 0.0000857, 0.00009737, 0.00009379, 0.00009407, 0.00009376, 0.00009904, 0.00010648, 0.00011187, 0.00010652, 0.00010911, 0.00010648, 0.00011525, 0.00011784, 0.00012511, 0.00013699, 0.00012817, 0.00014273, 0.00014226, 0.00014523, 0.00015219, 0.00015539, 0.00014564, 0.00015427, 0.00013868, 0.00014966, 0.00015355, 0.00016438, 0.00016471, 0.00016654, 0.0001928, 0.00018164, 0.00018786, 0.00020113, 0.00017872, 0.00019747, 0.00019327, 0.00019815, 0.00021518, 0.0002017, 0.00019843, 0.00022348, 0.00023552, 0.00025533, 0.00023887, 0.0002415, 0.00024799, 0.00024486, 0.00023786, 0.00023958, 0.0002336, 0.00027222, 0.0002468, 0.0002709, 0.00026972, 0.00029104, 0.00028203, 0.00027988, 0.00030318, 0.00029028, 0.00030375, 0.00029859, 0.00030314, 0.00028795, 0.00028249, 0.00030773, 0.00031414, 0.00032653, 0.00031658, 0.00033624, 0.00033422, 0.00034574, 0.00033151, 0.00033756, 0.00033123, 0.00037483, 0.00035685, 0.00036181, 0.0003705, 0.00035675, 0.00035758, 0.00035565, 0.00036023, 0.0003806, 0.00038147, 0.00038071, 0.00041048, 0.00039893, 0.00039449, 0.00040161, 0.00042135, 0.00041524, 0.0004308, 0.00042253, 0.00045223, 0.00040856, 0.00042726, 0.00040004, 0.00039657, 0.00044891, 0.00042078, 0.00047648, 0.00046267, 0.00045702, 0.00043118, 0.00044739, 0.00049803, 0.00044179, 0.00047492, 0.00044302, 0.00048639, 0.00049219, 0.00052193, 0.00054701, 0.00052215, 0.00049969, 0.00051054, 0.00053304, 0.00051857, 0.00053771]"

# CC_di array
$ws.Range("H10").Value = "[3.5, 3.70833333, 5.58333333, 4.875, 3, 7.58333333, 4.16666667, 2.70833333, 3.5, 3.70833333, 2, -0.41666667, 30.75, 204.58333333, 877, 2649.33333333, 5710.66666667, 9345.375, 12556.29166667, 14892.54166667, 16482.5, 17624.75, 18592.875, 19596.5, 20713.04166667, 22083.95833333, 23729.16666667, 25663, 27915.54166667, 30383.41666667, 33021.875, 35663.16666667, 38128.79166667, 40056.58333333, 41383.25, 41940.33333333, 41629.5, 40619.54166667, 39009.91666667, 37137.70833333, 35170.29166667, 33330.75, 31711.375, 30359.45833333, 29362.45833333, 28654.25, 28298.375, 28244.95833333, 28491.375, 29045.125, 29884.95833333, 31012.5, 32429, 34055.33333333, 35953.54166667, 37997.66666667, 40081.75, 42170.75, 44142.33333333, 45832.75, 47184.79166667, 48049.45833333, 48293.25, 47975.20833333, 47058, 45632.5, 43813.875, 41806.625, 39617.16666667, 37446.625, 35361, 33420.375, 31737.5, 30220.08333333, 28904.41666667, 27783.79166667, 26863.33333333, 26089.45833333, 25492.5, 25043.125, 24741.20833333, 24578.16666667, 24548.91666667, 24651.83333333, 24892, 25218.5, 25676.58333333, 26242.41666667, 26850.41666667, 27536.83333333, 28290.16666667, 29054.625, 29783.70833333, 30427.08333333, 30991.79166667, 31386, 31571.58333333, 31497.33333333, 31187.5, 30671.625, 29892.45833333, 28947.70833333, 27877.29166667, 26689.5, 25504.45833333, 24344.91666667, 23209.875, 22174.25, 21228.75, 20371.29166667, 19637.79166667, 18974.25, 18420.83333333, 17976.54166667, 17599.45833333, 17312.54166667, 17112, 16993.79166667, 16927.91666667, 16955.33333333, 17050.66666667, 17217.45833333, 17446.66666667, 17754.91666667, 18127.95833333, 18582.04166667, 19087.20833333, 19678.875, 20312.75, 20975.79166667, 21664.16666667, 22393.375, 23109.875, 23848.54166667, 24557.125, 25207.125, 25829.91666667, 26374.75, 26808.75, 27179.41666667, 27435.83333333, 27584.79166667, 27622.875, 27547.16666667, 27344.33333333, 27055.20833333, 26643.33333333, 26153.33333333, 25595.375, 24962.375, 24299.375, 23625.70833333, 22921.45833333, 22216.16666667, 21553, 20875.70833333, 20242.375, 19651.95833333, 19105.58333333, 18610.375, 18133.375, 17712.58333333, 17340.29166667, 16962.33333333, 16625.58333333, 16312.04166667, 15998.29166667, 15689.75, 15416.375, 15164.16666667, 14906.83333333, 14682.54166667, 14487.625, 14314.54166667, 14175.91666667, 14051.16666667, 13946.08333333, 13849.83333333, 13778.16666667, 13718.25, 13673.70833333, 13630.5, 13601.75, 13570.70833333, 13545.29166667, 13526.91666667, 13515.83333333, 13506.5, 13497.04166667, 13491.125, 13475.54166667, 13468.16666667, 13474.33333333, 13459.5, 13470.91666667, 13458.125, 13442.625, 13422.33333333, 13384.45833333, 13320.91666667, 13276.625, 13205.54166667, 13114.54166667, 13016.29166667, 12886.375, 12741.41666667, 12586.16666667, 12393.20833333, 12187.75, 11968.625, 11727.58333333, 11473.91666667, 11211.75, 10936.25, 10652.20833333, 10377.20833333, 10092.54166667, 9814.58333333, 9545.58333333, 9290.625, 9037.58333333, 8780.54166667, 8516.20833333, 8243.79166667, 7968.58333333, 7720.54166667, 7494.70833333, 7296.54166667, 7113.875, 6963.70833333, 6826.95833333, 6703, 6607.375, 6518.20833333, 6437.20833333, 6372.04166667, 6300.33333333, 6207.875, 6099.70833333, 5996.70833333, 5905.91666667, 5830.79166667, 5757.16666667, 5694.16666667, 5642.375, 5599.54166667, 5564.95833333, 5541.125, 5528.66666667, 5522.70833333, 5514.16666667, 5515.125, 5519.41666667, 5465.08333333, 5280.33333333, 4866.83333333]"

# CC_lower_wavelength_limit_for_spectra_fit - highlighted green, numeric format 0.00
$ws.Range("F6").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$ws.Range("H11").Value = 217
$ws.Range("H11").Interior.Color = 5296274
$ws.Range("H11").NumberFormat = "0.00"
$ws.Range("H11").HorizontalAlignment = -4131
$ws.Range("H11").VerticalAlignment = -4108

# CC_upper_wavelength_limit_for_spectra_fit - same new style, highlighted green
$ws.Range("H11").Copy()
$ws.Range("H12").PasteSpecial(-4122)
$ws.Range("H12").Value = 240
